# aligner_Error_Code.xlsx — "Add files via upload"
#
# The uploaded copy of the sheet gave the two label cells (C24 "얼라인 설정"
# and C25 "얼라인 실행") a second line repeating the label with no space
# (e.g. "얼라인 설정\n얼라인설정"), turned wrap-text on for those two cells,
# and left the selection sitting on the C24:D25 block that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 first, then row 24: the order the two new shared strings get
# (re)created in matters for matching shared-string table order downstream.
$ws.Range("C25").Value = "얼라인 실행`n얼라인실행"
$ws.Range("C24").Value = "얼라인 설정`n얼라인설정"

# Both label cells now hold two lines of text, so wrap them.
$ws.Range("C24:C25").WrapText = $true

# Leave the selection where the author left it after editing those cells,
# and scroll so row 23 is the first visible row.
$ws.Range("C24:D25").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
